$d = $word.ActiveDocument

# 1. Title change
$d.Content.Find.Execute(
    "Activity 1 - Why Investigate Inequity in Iowa?", $true, $false, $false, $false, $false,
    $true, 1, $false, "Activity 1 - What is Inequity?", 2)

# 2. Intro paragraph - expand on houseless community service description
$d.Content.Find.Execute(
    "Waypoint Services has tasked us to investigate a community service to help people who are houseless. Let",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Waypoint Services has tasked us to investigate a community service to help houseless people, specifically a program for people who have been incarcerated recently. Let",
    2)

# 3. "at one time" -> "at the same time"
$d.Content.Find.Execute(
    "you will all be able to work on that document at one time.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "you will all be able to work on that document at the same time.",
    2)

# 4. Remove sentence about Waypoint founding suspicions; "specifically" -> "explicitly"
$d.Content.Find.Execute(
    "please look for some anecdotal evidence. This type of evidence is primarily what Waypoint has founded their suspicions of inequity on. We are specifically looking for inequity in Iowa",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "please look for some anecdotal evidence. We are explicitly looking for inequity in Iowa",
    2)

# 5. "Start by talking among" -> "Begin by discussing with"; "key words" -> "keywords";
#    "natural disaster" -> "natural disasters"; drop trailing "Your goal is to"
$d.Content.Find.Execute(
    "Start by talking among your group how you will search: what key words will you use?",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Begin by discussing with your group how you will search: what keywords will you use?",
    2)

$d.Content.Find.Execute(
    "There have been two natural disaster in Iowa in the last 20 years, a flood and a derecho, both of which destroyed many homes. Consider these events when conducting your search. Your goal is to",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "There have been two natural disasters in Iowa in the last 20 years, a flood and a derecho, both of which destroyed many homes. Consider these events when conducting your search.",
    2)

# 6. "each sources" -> "each source"; parenthetical rewording
$d.Content.Find.Execute(
    "Make sure you are taking notes on each sources so you remember why you included it. (make sure to cite your sources, websites are fine in this case)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Make sure you are taking notes on each source so you remember why you included it. (Make sure to cite your sources; websites URLs are fine in this case.)",
    2)

# 7. Fix grammar in example-of-notes paragraph
$d.Content.Find.Execute(
    "They discuss how most of the that neighborhood are individual of color",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "They discuss how most of the people in that neighborhood are individuals of color.",
    2)

# 8. "10 to 20 examples" -> "5 to 10 examples"
$d.Content.Find.Execute(
    "Aim for at least 10 to 20 examples. You have a group for good reason.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Aim for at least 5 to 10 examples. You have a group for good reason.",
    2)

# 9. Presentation paragraph rewording
$d.Content.Find.Execute(
    "Create a presentation to be presented today on your findings. Use your examples to tell a story. You can use Powerpoint or other visualizing software but you are not required to. Make your presentation compelling and feel free to show off some of the most meaningful examples. All students must present some of your findings.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Create a presentation to be presented today on your findings (3min-5min). Use your examples to tell a story. You can use PowerPoint or other visualizing software. Make your presentation compelling, and feel free to show off some of the most meaningful examples. All students must present some of their findings.",
    2)

Write-Output $d.Content.Text
